$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D4").Value = 121740.8710009183
    $ws.Range("E4").Value = 0.001914870283549103
    $ws.Range("F4").Value = 0.1977593504539668
    $ws.Range("G4").Value = -1.066475700650384
    $ws.Range("H4").Value = 10.99633228977906
    $ws.Range("D5").Value = 122452.110392334
    $ws.Range("E5").Value = -0.00486428498044088
    $ws.Range("F5").Value = 0.2076432377418609
    $ws.Range("G5").Value = -0.5837990147931015
    $ws.Range("H5").Value = 7.054435814165112
    $ws.Range("D6").Value = 122983.8712583094
    $ws.Range("E6").Value = -0.01658695170572346
    $ws.Range("F6").Value = 0.243351302435544
    $ws.Range("G6").Value = -1.266485509600165
    $ws.Range("H6").Value = 11.32052294318789
    $ws.Range("D8").Value = 124668.1345378395
    $ws.Range("E8").Value = -0.03825809378727953
    $ws.Range("F8").Value = 0.2137051312979744
    $ws.Range("G8").Value = -0.8770110179930539
    $ws.Range("H8").Value = 6.970634387171684
    $ws.Range("D9").Value = 127031.0048926685
    $ws.Range("E9").Value = -0.06098362685232529
    $ws.Range("F9").Value = 0.3292842029811078
    $ws.Range("G9").Value = -1.612580553969292
    $ws.Range("H9").Value = 10.82731377842855
    $ws.Range("D10").Value = 128190.0269136586
    $ws.Range("E10").Value = -0.1009731354988247
    $ws.Range("F10").Value = 0.4245909280697951
    $ws.Range("G10").Value = -1.911376404514875
    $ws.Range("H10").Value = 10.07804588867599
    $ws.Range("D11").Value = 129973.5368707666
    $ws.Range("E11").Value = -0.1761469234126501
    $ws.Range("F11").Value = 0.7307554137921819
    $ws.Range("G11").Value = -2.587271510301969
    $ws.Range("H11").Value = 12.63140653012571
    $ws.Range("D13").Value = 120154.6788324331
    $ws.Range("E13").Value = 0.05215863132567373
    $ws.Range("F13").Value = 0.1182420371903864
    $ws.Range("G13").Value = -0.75575921568768
    $ws.Range("H13").Value = 7.109740548134341
    $ws.Range("D15").Value = 120175.6652923589
    $ws.Range("E15").Value = 0.0479361702042674
    $ws.Range("F15").Value = 0.1298370840022972
    $ws.Range("G15").Value = -0.7928914917491404
    $ws.Range("H15").Value = 10.8175258932869
    $ws.Range("D16").Value = 120154.9150019117
    $ws.Range("E16").Value = 0.1078258597277365
    $ws.Range("F16").Value = 0.1356155457900996
    $ws.Range("G16").Value = -0.2837711364968449
    $ws.Range("H16").Value = 4.341558542499821
    $ws.Range("D17").Value = 120122.5375086681
    $ws.Range("E17").Value = 0.08488983322559919
    $ws.Range("F17").Value = 0.1075199947981904
    $ws.Range("G17").Value = -0.9228967536253245
    $ws.Range("H17").Value = 8.495708435479653
    $ws.Range("D20").Value = 120968.5900089347
    $ws.Range("E20").Value = 0.03027171064179338
    $ws.Range("F20").Value = 0.1447422973101908
    $ws.Range("G20").Value = -0.2453646266519291
    $ws.Range("H20").Value = 6.073245055303178
